$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.240.09'
$ws.Range('D2').Style = $style
$ws.Range('E2').Value = '  -0.48%  '
$style = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.863.32'
$ws.Range('D3').Style = $style
$ws.Range('E3').Value = '  -1.13%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('E5').Value = '  -1.07%  '
$style = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.31'
$ws.Range('D6').Style = $style
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('E7').Value = '  +0.09%  '
$style = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07826'
$ws.Range('D8').Style = $style
$ws.Range('E8').Value = '  -3.04%  '
$style = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3108'
$ws.Range('D9').Style = $style
$ws.Range('E9').Value = '  -0.77%  '
$style = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.25'
$ws.Range('D10').Style = $style
$ws.Range('E10').Value = '  -4.06%  '
$style = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07989'
$ws.Range('D11').Style = $style
$ws.Range('E11').Value = '  -4.42%  '
$style = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.885.48'
$ws.Range('D12').Style = $style
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('B13').Value = 'Litecoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$style = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '93.63'
$ws.Range('D13').Style = $style
$ws.Range('E13').Value = '  +1.34%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$style = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.177'
$ws.Range('D14').Style = $style
$ws.Range('E14').Value = '  -1.31%  '
$style = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6950'
$ws.Range('D15').Style = $style
$ws.Range('E15').Value = '  -3.60%  '
$style = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.346'
$ws.Range('D16').Style = $style
$ws.Range('E16').Value = '  +1.07%  '
$style = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.388.26'
$ws.Range('D17').Style = $style
$ws.Range('E17').Value = '  +0.02%  '
$style = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008281'
$ws.Range('D18').Style = $style
$ws.Range('E18').Value = '  -2.17%  '
$style = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '252.11'
$ws.Range('D19').Style = $style
$ws.Range('E19').Value = '  +4.52%  '
$style = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.178.71'
$ws.Range('D20').Style = $style
$ws.Range('E20').Value = '  +2.57%  '
$ws.Range('E21').Value = '  -1.00%  '
$ws.Range('E22').Value = '  +0.05%  '
$style = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.539'
$ws.Range('D23').Style = $style
$ws.Range('E23').Value = '  -3.66%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('E25').Value = '  -2.20%  '
$style = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.992'
$ws.Range('D26').Style = $style
$ws.Range('E26').Value = '  -0.89%  '
$style = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '159.59'
$ws.Range('D27').Style = $style
$style = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.71'
$ws.Range('D28').Style = $style
$ws.Range('E28').Value = '  +0.68%  '
$style = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.498'
$ws.Range('D29').Style = $style
$ws.Range('E29').Value = '  -0.64%  '
$style = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.269'
$ws.Range('D30').Style = $style
$ws.Range('E30').Value = '  -1.58%  '
$style = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.267'
$ws.Range('D31').Style = $style
$ws.Range('E31').Value = '  -3.55%  '
$style = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.212'
$ws.Range('D32').Style = $style
$ws.Range('E32').Value = '  +0.60%  '
$style = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05259'
$ws.Range('D33').Style = $style
$ws.Range('E33').Value = '  -2.22%  '
$style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.893'
$ws.Range('D34').Style = $style
$ws.Range('E34').Value = '  -3.11%  '
$style = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7448'
$ws.Range('D35').Style = $style
$ws.Range('E35').Value = '  -0.71%  '
$style = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.156'
$ws.Range('D36').Style = $style
$ws.Range('E36').Value = '  -2.12%  '
$ws.Range('E37').Value = '  +0.11%  '
$ws.Range('E38').Value = '  -1.33%  '
$style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.247.50'
$ws.Range('D39').Style = $style
$ws.Range('E39').Value = '  -2.97%  '
$style = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.744'
$ws.Range('D40').Style = $style
$ws.Range('E40').Value = '  -0.11%  '
$style = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.260'
$ws.Range('D41').Style = $style
$ws.Range('E41').Value = '  -4.68%  '
$style = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9010'
$ws.Range('D42').Style = $style
$ws.Range('E42').Value = '  +0.90%  '
$style = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '111.12'
$ws.Range('D43').Style = $style
$ws.Range('E43').Value = '  +0.59%  '
$style = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '71.94'
$ws.Range('D44').Style = $style
$ws.Range('E44').Value = '  -2.17%  '
$style = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.001'
$ws.Range('D45').Style = $style
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$style = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000130'
$ws.Range('D46').Style = $style
$ws.Range('E46').Value = '  +0.61%  '
$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$style = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.067.59'
$ws.Range('D47').Style = $style
$ws.Range('E47').Value = '  +2.33%  '
$ws.Range('E48').Value = '  -0.28%  '
$style = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.786'
$ws.Range('D49').Style = $style
$ws.Range('E49').Value = '  -1.15%  '
$style = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.387'
$ws.Range('D50').Style = $style
$style = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.003'
$ws.Range('D51').Style = $style
$ws.Range('E51').Value = '  +0.02%  '
